$wb = $excel.ActiveWorkbook

# --- Sheet1 ("Days" timeline) ---
$ws1 = $wb.Worksheets.Item("Sheet1")

# Remove the first 10 days of data (rows 2-11); everything below shifts up,
# row/column numbering and dependent formulas (incl. cross-sheet refs) are
# recalculated automatically by Excel.
$ws1.Rows("2:11").Delete()

# Update the selected cell on Sheet1 to reflect where the author was working.
[void]$ws1.Range("B9").Select()

# --- Sheet2 ("Summary") ---
$ws2 = $wb.Worksheets.Item("Sheet2")

# Unit/row being restarted - zero out hours & projects, leave SubUnits as-is.
$ws2.Cells.Item(4, 2).Value = 0   # B4 Units
$ws2.Cells.Item(4, 4).Value = 0   # D4 Projects

# Add a note in column F for this row flagging it to revisit later.
$ws2.Cells.Item(4, 6).Value = "<- Come back to as needed"

# Update the selected cell on Sheet2 to reflect where the author was working.
[void]$ws2.Range("F20").Select()
